$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three data rows (2007, 2008, 2009) so that the
# remaining years (2010, 2011, 2012) shift up to become rows 2-4.
$ws.Range("A2:F4").Delete(-4162)

